$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: clone the number/cell formatting of a template cell onto a
#     target cell via copy / paste-special (formats only) so the engine
#     reuses the existing style index instead of minting a new one. ---
function Copy-Format($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial(-4122) | Out-Null
}

# ----------------------------------------------------------------------
# 1) Dates (column A) for the new rows 4-15. Style comes from row 2 (A2).
# ----------------------------------------------------------------------
$dates = @{
    4  = 45794
    5  = 45795
    6  = 45794
    7  = 45795
    8  = 45802
    9  = 45803
    10 = 45806
    11 = 45807
    12 = 45808
    13 = 45809
    14 = 45808
    15 = 45809
}
foreach ($r in 4..15) {
    Copy-Format "A2" "A$r"
    $ws.Range("A$r").Value2 = $dates[$r]
}

# ----------------------------------------------------------------------
# 2) Icon (column B) + name (column C) - first batch: the "saza ->
#    higedan -> g-dragon -> j-hope" rows, entered icon-column-first then
#    name-column-first (matches how the shared-string table grew).
# ----------------------------------------------------------------------
$firstBatchRows = 4, 5, 6, 7, 8, 9, 14, 15
$redIcon = "🔴"
$blueIcon = "🔵"
$icons = @{
    4  = $redIcon
    5  = $redIcon
    6  = $blueIcon
    7  = $blueIcon
    8  = $redIcon
    9  = $redIcon
    14 = $redIcon
    15 = $redIcon
}
$names = @{
    4  = "サザンオールスターズ"
    5  = "サザンオールスターズ"
    6  = "髭男dism"
    7  = "髭男dism"
    8  = "G-DRAGON"
    9  = "G-DRAGON"
    14 = "J-HOPE"
    15 = "J-HOPE"
}

foreach ($r in $firstBatchRows) {
    Copy-Format "B2" "B$r"
    $ws.Range("B$r").Value = $icons[$r]
}
foreach ($r in $firstBatchRows) {
    $ws.Range("C$r").Value = $names[$r]
}

# ----------------------------------------------------------------------
# 3) Icon + name - second batch: the PokeGoFes rows (10-13), inserted
#    after the first batch.
# ----------------------------------------------------------------------
$secondBatchRows = 10, 11, 12, 13
$star = "★"

foreach ($r in $secondBatchRows) {
    Copy-Format "B2" "B$r"
    $ws.Range("B$r").Value = $star
}
foreach ($r in $secondBatchRows) {
    $ws.Range("C$r").Value = "ポケGofes"
}

# ----------------------------------------------------------------------
# 4) Header row: B1 stays "icon" (value unchanged, just a shared-string
#    reindex happens naturally). Finally fix the SEVENTEEN LIVE rows
#    (drop the stray leading tab) - done last, matching the order in
#    which the shared-string table was rebuilt.
# ----------------------------------------------------------------------
$ws.Range("B1").Value = "icon"
$ws.Range("C2").Value = "SEVENTEEN LIVE"
$ws.Range("C3").Value = "SEVENTEEN LIVE"

# ----------------------------------------------------------------------
# 5) Selection, matching the saved workbook view.
# ----------------------------------------------------------------------
$ws.Range("C3").Select() | Out-Null
